$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for rows 2 through 31.
# Update every one of these cells from 45588 (2024-10-23) to 45589 (2024-10-24),
# keeping the existing number formatting/style untouched.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45588) {
        $cell.Value2 = 45589
    }
}
